$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily date entries appended to the "Diaria" (daily) table,
# following the same pattern as the existing rows (A: date, B: 10000, D: 0).
$newDates = @("21-09-2021", "22-09-2021", "23-09-2021", "28-09-2021", "30-09-2021")

$startRow = 25
for ($i = 0; $i -lt $newDates.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newDates[$i]
    $ws.Cells.Item($row, 2).Value = 10000
    $ws.Cells.Item($row, 4).Value = 0
}

$wb.Save()
